{"js": "const body = context.document.body;\nconst results = body.search(\"Scenario registrovanja korisnika\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const range = results.items[0];\n  range.insertText(\"Scenario kreiranje grupe od strane admina\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Replacement.ClearFormatting()\n$find.Find.Execute(\n    \"Scenario registrovanja korisnika\",  # FindText\n    $false,                               # MatchCase\n    $false,                               # MatchWholeWord\n    $false,                               # MatchWildcards\n    $false,                               # MatchSoundsLike\n    $false,                               # MatchAllWordForms\n    $true,                                # Forward\n    1,                                    # Wrap (wdFindContinue)\n    $false,                               # Format\n    \"Scenario kreiranje grupe od strane admina\",  # ReplaceWith\n    2                                     # Replace (wdReplaceAll)\n)\n"}
